$d = $word.ActiveDocument

# Update the date heading at the top of the document
$d.Content.Find.Execute("2023-09-04 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-05 Tuesday", 1) | Out-Null

# Update each arithmetic-problem cell in the table, position by position.
# wdReplaceOne (1) is used (not wdReplaceAll) so the replace is confined to the
# current cell even when the same expression text recurs elsewhere in the table.
$t = $d.Tables.Item(1)

$cell = $t.Cell(1,1)
$cell.Range.Find.Execute("40-9=", $true, $false, $false, $false, $false, $true, 1, $false, "34+19=", 1) | Out-Null
$cell = $t.Cell(1,2)
$cell.Range.Find.Execute("10+20=", $true, $false, $false, $false, $false, $true, 1, $false, "29+55=", 1) | Out-Null
$cell = $t.Cell(1,3)
$cell.Range.Find.Execute("83-81=", $true, $false, $false, $false, $false, $true, 1, $false, "24+40=", 1) | Out-Null
$cell = $t.Cell(1,4)
$cell.Range.Find.Execute("50-44=", $true, $false, $false, $false, $false, $true, 1, $false, "8+14=", 1) | Out-Null
$cell = $t.Cell(1,5)
$cell.Range.Find.Execute("17-11=", $true, $false, $false, $false, $false, $true, 1, $false, "32-18=", 1) | Out-Null

$cell = $t.Cell(2,1)
$cell.Range.Find.Execute("45+27=", $true, $false, $false, $false, $false, $true, 1, $false, "38+28=", 1) | Out-Null
$cell = $t.Cell(2,2)
$cell.Range.Find.Execute("37+53=", $true, $false, $false, $false, $false, $true, 1, $false, "33+51=", 1) | Out-Null
$cell = $t.Cell(2,3)
$cell.Range.Find.Execute("94-63=", $true, $false, $false, $false, $false, $true, 1, $false, "51-23=", 1) | Out-Null
$cell = $t.Cell(2,4)
$cell.Range.Find.Execute("94-28=", $true, $false, $false, $false, $false, $true, 1, $false, "82-42=", 1) | Out-Null
$cell = $t.Cell(2,5)
$cell.Range.Find.Execute("69-37=", $true, $false, $false, $false, $false, $true, 1, $false, "47+0=", 1) | Out-Null

$cell = $t.Cell(3,1)
$cell.Range.Find.Execute("99-12=", $true, $false, $false, $false, $false, $true, 1, $false, "18-13=", 1) | Out-Null
$cell = $t.Cell(3,2)
$cell.Range.Find.Execute("40-21=", $true, $false, $false, $false, $false, $true, 1, $false, "44+2=", 1) | Out-Null
$cell = $t.Cell(3,3)
$cell.Range.Find.Execute("69-29=", $true, $false, $false, $false, $false, $true, 1, $false, "36-31=", 1) | Out-Null
$cell = $t.Cell(3,4)
$cell.Range.Find.Execute("34-25=", $true, $false, $false, $false, $false, $true, 1, $false, "61-6=", 1) | Out-Null
$cell = $t.Cell(3,5)
$cell.Range.Find.Execute("43-14=", $true, $false, $false, $false, $false, $true, 1, $false, "80+8=", 1) | Out-Null

$cell = $t.Cell(4,1)
$cell.Range.Find.Execute("59-39=", $true, $false, $false, $false, $false, $true, 1, $false, "5+20=", 1) | Out-Null
$cell = $t.Cell(4,2)
$cell.Range.Find.Execute("34-21=", $true, $false, $false, $false, $false, $true, 1, $false, "98-44=", 1) | Out-Null
$cell = $t.Cell(4,3)
$cell.Range.Find.Execute("62+23=", $true, $false, $false, $false, $false, $true, 1, $false, "71-62=", 1) | Out-Null
$cell = $t.Cell(4,4)
$cell.Range.Find.Execute("28+57=", $true, $false, $false, $false, $false, $true, 1, $false, "9+90=", 1) | Out-Null
$cell = $t.Cell(4,5)
$cell.Range.Find.Execute("31-26=", $true, $false, $false, $false, $false, $true, 1, $false, "39+22=", 1) | Out-Null

$cell = $t.Cell(5,1)
$cell.Range.Find.Execute("6+51=", $true, $false, $false, $false, $false, $true, 1, $false, "8+51=", 1) | Out-Null
$cell = $t.Cell(5,2)
$cell.Range.Find.Execute("35-31=", $true, $false, $false, $false, $false, $true, 1, $false, "31+28=", 1) | Out-Null
$cell = $t.Cell(5,3)
$cell.Range.Find.Execute("61-57=", $true, $false, $false, $false, $false, $true, 1, $false, "54+12=", 1) | Out-Null
$cell = $t.Cell(5,4)
$cell.Range.Find.Execute("12-9=", $true, $false, $false, $false, $false, $true, 1, $false, "38+8=", 1) | Out-Null
$cell = $t.Cell(5,5)
$cell.Range.Find.Execute("22+13=", $true, $false, $false, $false, $false, $true, 1, $false, "71-16=", 1) | Out-Null

$cell = $t.Cell(6,1)
$cell.Range.Find.Execute("86-44=", $true, $false, $false, $false, $false, $true, 1, $false, "39+49=", 1) | Out-Null
$cell = $t.Cell(6,2)
$cell.Range.Find.Execute("96-37=", $true, $false, $false, $false, $false, $true, 1, $false, "33+3=", 1) | Out-Null
$cell = $t.Cell(6,3)
$cell.Range.Find.Execute("52-1=", $true, $false, $false, $false, $false, $true, 1, $false, "15+73=", 1) | Out-Null
$cell = $t.Cell(6,4)
$cell.Range.Find.Execute("23+16=", $true, $false, $false, $false, $false, $true, 1, $false, "60+9=", 1) | Out-Null
$cell = $t.Cell(6,5)
$cell.Range.Find.Execute("11+57=", $true, $false, $false, $false, $false, $true, 1, $false, "51-24=", 1) | Out-Null

$cell = $t.Cell(7,1)
$cell.Range.Find.Execute("46+0=", $true, $false, $false, $false, $false, $true, 1, $false, "69+24=", 1) | Out-Null
$cell = $t.Cell(7,2)
$cell.Range.Find.Execute("76-42=", $true, $false, $false, $false, $false, $true, 1, $false, "58-30=", 1) | Out-Null
$cell = $t.Cell(7,3)
$cell.Range.Find.Execute("31+22=", $true, $false, $false, $false, $false, $true, 1, $false, "62-40=", 1) | Out-Null
$cell = $t.Cell(7,4)
$cell.Range.Find.Execute("65+16=", $true, $false, $false, $false, $false, $true, 1, $false, "32+35=", 1) | Out-Null
$cell = $t.Cell(7,5)
$cell.Range.Find.Execute("25+52=", $true, $false, $false, $false, $false, $true, 1, $false, "89-81=", 1) | Out-Null

$cell = $t.Cell(8,1)
$cell.Range.Find.Execute("13+68=", $true, $false, $false, $false, $false, $true, 1, $false, "72-43=", 1) | Out-Null
$cell = $t.Cell(8,2)
$cell.Range.Find.Execute("64-62=", $true, $false, $false, $false, $false, $true, 1, $false, "48+30=", 1) | Out-Null
$cell = $t.Cell(8,3)
$cell.Range.Find.Execute("0+65=", $true, $false, $false, $false, $false, $true, 1, $false, "44-36=", 1) | Out-Null
$cell = $t.Cell(8,4)
$cell.Range.Find.Execute("92-10=", $true, $false, $false, $false, $false, $true, 1, $false, "21+34=", 1) | Out-Null
$cell = $t.Cell(8,5)
$cell.Range.Find.Execute("99-58=", $true, $false, $false, $false, $false, $true, 1, $false, "59+25=", 1) | Out-Null

$cell = $t.Cell(9,1)
$cell.Range.Find.Execute("71+10=", $true, $false, $false, $false, $false, $true, 1, $false, "81-31=", 1) | Out-Null
$cell = $t.Cell(9,2)
$cell.Range.Find.Execute("39-17=", $true, $false, $false, $false, $false, $true, 1, $false, "66-15=", 1) | Out-Null
$cell = $t.Cell(9,3)
$cell.Range.Find.Execute("30+39=", $true, $false, $false, $false, $false, $true, 1, $false, "55+5=", 1) | Out-Null
$cell = $t.Cell(9,4)
$cell.Range.Find.Execute("62+35=", $true, $false, $false, $false, $false, $true, 1, $false, "54-39=", 1) | Out-Null
$cell = $t.Cell(9,5)
$cell.Range.Find.Execute("63-47=", $true, $false, $false, $false, $false, $true, 1, $false, "72+18=", 1) | Out-Null

$cell = $t.Cell(10,1)
$cell.Range.Find.Execute("66-21=", $true, $false, $false, $false, $false, $true, 1, $false, "73-69=", 1) | Out-Null
$cell = $t.Cell(10,2)
$cell.Range.Find.Execute("76-8=", $true, $false, $false, $false, $false, $true, 1, $false, "14+61=", 1) | Out-Null
$cell = $t.Cell(10,3)
$cell.Range.Find.Execute("34-13=", $true, $false, $false, $false, $false, $true, 1, $false, "25+27=", 1) | Out-Null
$cell = $t.Cell(10,4)
$cell.Range.Find.Execute("62+18=", $true, $false, $false, $false, $false, $true, 1, $false, "56-27=", 1) | Out-Null
$cell = $t.Cell(10,5)
$cell.Range.Find.Execute("40+32=", $true, $false, $false, $false, $false, $true, 1, $false, "68-25=", 1) | Out-Null

$cell = $t.Cell(11,1)
$cell.Range.Find.Execute("50+35=", $true, $false, $false, $false, $false, $true, 1, $false, "44+0=", 1) | Out-Null
$cell = $t.Cell(11,2)
$cell.Range.Find.Execute("58-52=", $true, $false, $false, $false, $false, $true, 1, $false, "29+13=", 1) | Out-Null
$cell = $t.Cell(11,3)
$cell.Range.Find.Execute("80-3=", $true, $false, $false, $false, $false, $true, 1, $false, "29+33=", 1) | Out-Null
$cell = $t.Cell(11,4)
$cell.Range.Find.Execute("77-11=", $true, $false, $false, $false, $false, $true, 1, $false, "80-16=", 1) | Out-Null
$cell = $t.Cell(11,5)
$cell.Range.Find.Execute("17+38=", $true, $false, $false, $false, $false, $true, 1, $false, "64-4=", 1) | Out-Null

$cell = $t.Cell(12,1)
$cell.Range.Find.Execute("46-19=", $true, $false, $false, $false, $false, $true, 1, $false, "11+56=", 1) | Out-Null
$cell = $t.Cell(12,2)
$cell.Range.Find.Execute("58+34=", $true, $false, $false, $false, $false, $true, 1, $false, "53-10=", 1) | Out-Null
$cell = $t.Cell(12,3)
$cell.Range.Find.Execute("84-53=", $true, $false, $false, $false, $false, $true, 1, $false, "0+5=", 1) | Out-Null
$cell = $t.Cell(12,4)
$cell.Range.Find.Execute("72-0=", $true, $false, $false, $false, $false, $true, 1, $false, "10+43=", 1) | Out-Null
$cell = $t.Cell(12,5)
$cell.Range.Find.Execute("41-25=", $true, $false, $false, $false, $false, $true, 1, $false, "27+50=", 1) | Out-Null

$cell = $t.Cell(13,1)
$cell.Range.Find.Execute("95-77=", $true, $false, $false, $false, $false, $true, 1, $false, "53-15=", 1) | Out-Null
$cell = $t.Cell(13,2)
$cell.Range.Find.Execute("30+69=", $true, $false, $false, $false, $false, $true, 1, $false, "24+52=", 1) | Out-Null
$cell = $t.Cell(13,3)
$cell.Range.Find.Execute("44+13=", $true, $false, $false, $false, $false, $true, 1, $false, "5+53=", 1) | Out-Null
$cell = $t.Cell(13,4)
$cell.Range.Find.Execute("17+29=", $true, $false, $false, $false, $false, $true, 1, $false, "80-36=", 1) | Out-Null
$cell = $t.Cell(13,5)
$cell.Range.Find.Execute("64+33=", $true, $false, $false, $false, $false, $true, 1, $false, "17+69=", 1) | Out-Null

$cell = $t.Cell(14,1)
$cell.Range.Find.Execute("86-37=", $true, $false, $false, $false, $false, $true, 1, $false, "4+46=", 1) | Out-Null
$cell = $t.Cell(14,2)
$cell.Range.Find.Execute("74-2=", $true, $false, $false, $false, $false, $true, 1, $false, "17+37=", 1) | Out-Null
$cell = $t.Cell(14,3)
$cell.Range.Find.Execute("28+57=", $true, $false, $false, $false, $false, $true, 1, $false, "36+0=", 1) | Out-Null
$cell = $t.Cell(14,4)
$cell.Range.Find.Execute("93-15=", $true, $false, $false, $false, $false, $true, 1, $false, "89-58=", 1) | Out-Null
$cell = $t.Cell(14,5)
$cell.Range.Find.Execute("74-72=", $true, $false, $false, $false, $false, $true, 1, $false, "74-33=", 1) | Out-Null

$cell = $t.Cell(15,1)
$cell.Range.Find.Execute("60+3=", $true, $false, $false, $false, $false, $true, 1, $false, "44+37=", 1) | Out-Null
$cell = $t.Cell(15,2)
$cell.Range.Find.Execute("34-21=", $true, $false, $false, $false, $false, $true, 1, $false, "58+7=", 1) | Out-Null
$cell = $t.Cell(15,3)
$cell.Range.Find.Execute("0+70=", $true, $false, $false, $false, $false, $true, 1, $false, "24+71=", 1) | Out-Null
$cell = $t.Cell(15,4)
$cell.Range.Find.Execute("34+6=", $true, $false, $false, $false, $false, $true, 1, $false, "5+44=", 1) | Out-Null
$cell = $t.Cell(15,5)
$cell.Range.Find.Execute("17-10=", $true, $false, $false, $false, $false, $true, 1, $false, "60-18=", 1) | Out-Null

$cell = $t.Cell(16,1)
$cell.Range.Find.Execute("76-12=", $true, $false, $false, $false, $false, $true, 1, $false, "32+2=", 1) | Out-Null
$cell = $t.Cell(16,2)
$cell.Range.Find.Execute("15+50=", $true, $false, $false, $false, $false, $true, 1, $false, "2+23=", 1) | Out-Null
$cell = $t.Cell(16,3)
$cell.Range.Find.Execute("64-42=", $true, $false, $false, $false, $false, $true, 1, $false, "26+40=", 1) | Out-Null
$cell = $t.Cell(16,4)
$cell.Range.Find.Execute("66+9=", $true, $false, $false, $false, $false, $true, 1, $false, "97-35=", 1) | Out-Null
$cell = $t.Cell(16,5)
$cell.Range.Find.Execute("26+38=", $true, $false, $false, $false, $false, $true, 1, $false, "2+40=", 1) | Out-Null

$cell = $t.Cell(17,1)
$cell.Range.Find.Execute("38-13=", $true, $false, $false, $false, $false, $true, 1, $false, "19-2=", 1) | Out-Null
$cell = $t.Cell(17,2)
$cell.Range.Find.Execute("94-43=", $true, $false, $false, $false, $false, $true, 1, $false, "78+16=", 1) | Out-Null
$cell = $t.Cell(17,3)
$cell.Range.Find.Execute("90-28=", $true, $false, $false, $false, $false, $true, 1, $false, "54-47=", 1) | Out-Null
$cell = $t.Cell(17,4)
$cell.Range.Find.Execute("93-65=", $true, $false, $false, $false, $false, $true, 1, $false, "56+35=", 1) | Out-Null
$cell = $t.Cell(17,5)
$cell.Range.Find.Execute("85-61=", $true, $false, $false, $false, $false, $true, 1, $false, "30+64=", 1) | Out-Null

$cell = $t.Cell(18,1)
$cell.Range.Find.Execute("87-13=", $true, $false, $false, $false, $false, $true, 1, $false, "30+42=", 1) | Out-Null
$cell = $t.Cell(18,2)
$cell.Range.Find.Execute("46+44=", $true, $false, $false, $false, $false, $true, 1, $false, "50-37=", 1) | Out-Null
$cell = $t.Cell(18,3)
$cell.Range.Find.Execute("36+18=", $true, $false, $false, $false, $false, $true, 1, $false, "73-65=", 1) | Out-Null
$cell = $t.Cell(18,4)
$cell.Range.Find.Execute("60+33=", $true, $false, $false, $false, $false, $true, 1, $false, "87-85=", 1) | Out-Null
$cell = $t.Cell(18,5)
$cell.Range.Find.Execute("64+9=", $true, $false, $false, $false, $false, $true, 1, $false, "46+32=", 1) | Out-Null

$cell = $t.Cell(19,1)
$cell.Range.Find.Execute("4+8=", $true, $false, $false, $false, $false, $true, 1, $false, "49-32=", 1) | Out-Null
$cell = $t.Cell(19,2)
$cell.Range.Find.Execute("87-4=", $true, $false, $false, $false, $false, $true, 1, $false, "27+50=", 1) | Out-Null
$cell = $t.Cell(19,3)
$cell.Range.Find.Execute("47+38=", $true, $false, $false, $false, $false, $true, 1, $false, "47+13=", 1) | Out-Null
$cell = $t.Cell(19,4)
$cell.Range.Find.Execute("59-7=", $true, $false, $false, $false, $false, $true, 1, $false, "16+22=", 1) | Out-Null
$cell = $t.Cell(19,5)
$cell.Range.Find.Execute("89-12=", $true, $false, $false, $false, $false, $true, 1, $false, "4+87=", 1) | Out-Null

$cell = $t.Cell(20,1)
$cell.Range.Find.Execute("6+18=", $true, $false, $false, $false, $false, $true, 1, $false, "26+60=", 1) | Out-Null
$cell = $t.Cell(20,2)
$cell.Range.Find.Execute("54-24=", $true, $false, $false, $false, $false, $true, 1, $false, "54-36=", 1) | Out-Null
$cell = $t.Cell(20,3)
$cell.Range.Find.Execute("71-18=", $true, $false, $false, $false, $false, $true, 1, $false, "14+11=", 1) | Out-Null
$cell = $t.Cell(20,4)
$cell.Range.Find.Execute("95-94=", $true, $false, $false, $false, $false, $true, 1, $false, "74-55=", 1) | Out-Null
$cell = $t.Cell(20,5)
$cell.Range.Find.Execute("76-54=", $true, $false, $false, $false, $false, $true, 1, $false, "0+90=", 1) | Out-Null
